# river update May 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("F2").Value = 0.897662974082117
$ws.Range("G2").Value = 0.0238095238095238
$ws.Range("H2").Value = 0.928571428571429
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = -0.802747252747253
$ws.Range("L2").Value = -2.12226793813417
$ws.Range("M2").Value = 0.101241745549664
$ws.Range("N2").Value = -21.9930754177329
$ws.Range("P2").Value = "Likely improving"

# Row 3 updates
$ws.Range("F3").Value = 0.1561846824104
$ws.Range("G3").Value = 0.0105263157894737
$ws.Range("H3").Value = 0.747368421052632
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 0.0454229257864196
$ws.Range("L3").Value = -0.0423246827733439
$ws.Range("M3").Value = 0.209048047845462
$ws.Range("N3").Value = 1.46525567052967

# Row 4 updates
$ws.Range("F4").Value = 0.08827525718496421
$ws.Range("G4").Value = 0.0068493150684931
$ws.Range("H4").Value = 0.712328767123288
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2.775
$ws.Range("K4").Value = 0.0265119398274774
$ws.Range("L4").Value = -0.0034724443135832
$ws.Range("M4").Value = 0.0971633583646566
$ws.Range("N4").Value = 0.955385219008194
$ws.Range("P4").Value = "Very unlikely improving"
